# Swap the presentation's applied colour theme from the custom "Integral"
# palette to the default "Office Theme" palette (the 12-slot theme colour
# scheme carried on the slide master's theme).
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$cs = $master.Theme.ThemeColorScheme

# MsoThemeColorSchemeIndex order: 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2
# 7=accent3 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink.
# RGB() is not available in this host, so the values below are the
# equivalent 0xBBGGRR long values for the standard Office theme colours.
$cs.Item(1).RGB  = 0x000000   # dk1      000000
$cs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$cs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$cs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink 954F72
